$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '60.362.91'
$ws.Range('E2').Value = '  +1.93%  '

$ws.Range('D3').Value = '2.538.92'
$ws.Range('E3').Value = '  +2.23%  '

$ws.Range('E4').Value = '  -0.01%  '

$ws.Range('D5').Value = '541.44'
$ws.Range('E5').Value = '  +0.55%  '

$ws.Range('D6').Value = '145.47'
$ws.Range('E6').Value = '  -0.20%  '

$ws.Range('E7').Value = '  -0.27%  '

$ws.Range('D8').Value = '0.574'
$ws.Range('E8').Value = '  +0.18%  '

$ws.Range('D9').Value = '2.571.32'
$ws.Range('E9').Value = '  +2.39%  '

$ws.Range('E10').Value = '  +1.08%  '

$ws.Range('E11').Value = '  +1.60%  '

$ws.Range('E12').Value = '  -1.16%  '

$ws.Range('D13').Value = '0.364'
$ws.Range('E13').Value = '  +1.74%  '

$ws.Range('D14').Value = '2.982.48'
$ws.Range('E14').Value = '  +1.87%  '

$ws.Range('D15').Value = '24.26'
$ws.Range('E15').Value = '  +0.47%  '

$ws.Range('D16').Value = '60.281.73'
$ws.Range('E16').Value = '  +1.97%  '

$ws.Range('E17').Value = '  +3.46%  '

$ws.Range('D18').Value = '2.550.13'
$ws.Range('E18').Value = '  +1.54%  '

$ws.Range('D19').Value = '11.38'
$ws.Range('E19').Value = '  -0.94%  '

$ws.Range('D20').Value = '4.37'
$ws.Range('E20').Value = '  +0.38%  '

$ws.Range('D21').Value = '329.10'
$ws.Range('E21').Value = '  +0.89%  '

$ws.Range('D22').Value = '1.00'
$ws.Range('E22').Value = '  +0.22%  '

$ws.Range('E23').Value = '  +2.31%  '

$ws.Range('D24').Value = '62.89'
$ws.Range('E24').Value = '  +2.87%  '

$ws.Range('D25').Value = '0.441'
$ws.Range('E25').Value = '  -0.80%  '

$ws.Range('E26').Value = '  +3.27%  '

$ws.Range('E27').Value = '  -0.42%  '

$ws.Range('E28').Value = '  +1.87%  '

$ws.Range('D29').Value = '7.17'
$ws.Range('E29').Value = '  +1.03%  '

$ws.Range('D30').Value = '0.0₃0801'
$ws.Range('E30').Value = '  +2.28%  '

$ws.Range('E31').Value = '  -0.06%  '

$ws.Range('E32').Value = '  -5.06%  '

$ws.Range('B33').Value = 'Monero'
$ws.Range('C33').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D33').Value = '162.67'
$ws.Range('E33').Value = '  +2.76%  '

$ws.Range('B34').Value = 'ImmutableX'
$ws.Range('C34').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D34').Value = '1.50'
$ws.Range('E34').Value = '  +5.63%  '

$ws.Range('E35').Value = '  +0.05%  '

$ws.Range('D36').Value = '18.84'
$ws.Range('E36').Value = '  +1.20%  '

$ws.Range('D37').Value = '4.52'
$ws.Range('E37').Value = '  +1.06%  '

$ws.Range('E38').Value = '  +0.06%  '

$ws.Range('D39').Value = '5.71'
$ws.Range('E39').Value = '  -2.88%  '

$ws.Range('D40').Value = '37.26'
$ws.Range('E40').Value = '  +1.27%  '

$ws.Range('D41').Value = '305.16'
$ws.Range('E41').Value = '  -2.80%  '

$ws.Range('D42').Value = '0.843'
$ws.Range('E42').Value = '  +1.46%  '

$ws.Range('D43').Value = '3.76'
$ws.Range('E43').Value = '  +0.42%  '

$ws.Range('B44').Value = 'Mantle'
$ws.Range('C44').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D44').Value = '0.609'
$ws.Range('E44').Value = '  +1.76%  '

$ws.Range('B45').Value = 'FirstDigitalUSD'
$ws.Range('C45').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D45').Value = '0.991'
$ws.Range('E45').Value = '  -0.38%  '

$ws.Range('E46').Value = '  +0.90%  '

$ws.Range('D47').Value = '19.14'
$ws.Range('E47').Value = '  +2.96%  '

$ws.Range('D48').Value = '0.0941'
$ws.Range('E48').Value = '  +1.04%  '

$ws.Range('D49').Value = '124.79'
$ws.Range('E49').Value = '  -0.42%  '

$ws.Range('D50').Value = '0.0525'
$ws.Range('E50').Value = '  -0.37%  '

$ws.Range('D51').Value = '0.0231'
$ws.Range('E51').Value = '  +0.22%  '
